$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(1003, "agrahari",   "Agr123456",     "agr@gmail.com"),
    @(1004, "agraalkll",  "jlkjlGu987978", "agr1@test.com"),
    @(1006, "agrahari78", "Agshui87987",   "agr3@gmail.com")
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    for ($c = 5; $c -le 15; $c++) {
        $ws.Cells.Item($r, $c).Value = $true
    }
    $ws.Cells.Item($r, 16).Value = $false
    $r++
}
